$wb = $excel.ActiveWorkbook

# Sheet 1: "peliculas o documentales" - remove the "Guerrilla del Oro" row (row 2),
# shifting the rows below it up.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("2:2").Delete()

# Sheet 2: "series" - remove the "My Brilliant Friend" row (row 2),
# shifting the rows below it up.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("2:2").Delete()
